$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings remain text (preserve exact formatting)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply updated values
$ws.Range("D2").Value = '41.861.79'
$ws.Range("E2").Value = '  -0.51%  '
$ws.Range("D3").Value = '2.273.16'
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '305.63'
$ws.Range("E5").Value = '  +1.13%  '
$ws.Range("D6").Value = '93.18'
$ws.Range("E6").Value = '  +0.21%  '
$ws.Range("E7").Value = '  -0.52%  '
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("E9").Value = '  +0.44%  '
$ws.Range("D10").Value = '32.70'
$ws.Range("E10").Value = '  -0.19%  '
$ws.Range("E11").Value = '  -0.31%  '
$ws.Range("E12").Value = '  -1.94%  '
$ws.Range("D13").Value = '6.69'
$ws.Range("E13").Value = '  -0.06%  '
$ws.Range("D14").Value = '2.624.28'
$ws.Range("E14").Value = '  +0.46%  '
$ws.Range("D15").Value = '14.37'
$ws.Range("E15").Value = '  +1.43%  '
$ws.Range("D16").Value = '2.273.77'
$ws.Range("E16").Value = '  +0.46%  '
$ws.Range("D17").Value = '0.784'
$ws.Range("E17").Value = '  +3.34%  '
$ws.Range("D18").Value = '41.783.99'
$ws.Range("E18").Value = '  -0.42%  '
$ws.Range("D19").Value = '12.79'
$ws.Range("E19").Value = '  +4.74%  '
$ws.Range("E20").Value = '  +0.81%  '
$ws.Range("D21").Value = '5.98'
$ws.Range("E21").Value = '  +0.28%  '
$ws.Range("D22").Value = '68.03'
$ws.Range("E22").Value = '  +1.06%  '
$ws.Range("D23").Value = '244.07'
$ws.Range("E23").Value = '  +1.04%  '
$ws.Range("E24").Value = '  +0.30%  '
$ws.Range("E25").Value = '  +1.98%  '
$ws.Range("E26").Value = '  +0.06%  '
$ws.Range("D27").Value = '24.03'
$ws.Range("E27").Value = '  +0.05%  '
$ws.Range("E28").Value = '  -0.15%  '
$ws.Range("E29").Value = '  -9.84%  '
$ws.Range("D30").Value = '34.73'
$ws.Range("E30").Value = '  +1.32%  '
$ws.Range("D31").Value = '159.53'
$ws.Range("E31").Value = '  +0.12%  '
$ws.Range("D32").Value = '5.37'
$ws.Range("E32").Value = '  +3.93%  '
$ws.Range("E33").Value = '  +0.00%  '
$ws.Range("E34").Value = '  +0.00%  '
$ws.Range("E35").Value = '  -0.89%  '
$ws.Range("D36").Value = '17.21'
$ws.Range("E36").Value = '  +3.54%  '
$ws.Range("E37").Value = '  -1.39%  '
$ws.Range("E38").Value = '  +0.42%  '
$ws.Range("E39").Value = '  +0.46%  '
$ws.Range("D40").Value = '1.80'
$ws.Range("E40").Value = '  -0.77%  '
$ws.Range("E41").Value = '  -0.26%  '
$ws.Range("B42").Value = 'EnergySwap'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D42").Value = '19.64'
$ws.Range("E42").Value = '  -2.12%  '
$ws.Range("B43").Value = 'Maker'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D43").Value = '2.006.34'
$ws.Range("E43").Value = '  -2.31%  '
$ws.Range("E44").Value = '  +12.05%  '
$ws.Range("E45").Value = '  +0.85%  '
$ws.Range("D46").Value = '10.26'
$ws.Range("E46").Value = '  +1.27%  '
$ws.Range("E47").Value = '  +0.29%  '
$ws.Range("D48").Value = '53.57'
$ws.Range("E48").Value = '  +3.01%  '
$ws.Range("E49").Value = '  +2.91%  '
$ws.Range("E50").Value = '  -1.16%  '
$ws.Range("D51").Value = '1.14'
$ws.Range("E51").Value = '  -0.11%  '
